$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Rename the second sheet from "F 0.5" to "customKcats" ---
$ws2.Name = "customKcats"

# --- Add newly curated kcat rows to the customKcats sheet (rows 33-37) ---
# Cell-write order mirrors the author's original entry order so the
# shared-string table grows in the same sequence as the source commit.
$ws2.Cells.Item(33, 1).Value = "Q6CGM4"
$ws2.Cells.Item(33, 2).Value = "YALI0A18062g"
$ws2.Cells.Item(33, 3).Value = "YALI0A18062g"
$ws2.Cells.Item(33, 4).Value = 0.584
$ws2.Cells.Item(33, 5).Value = "y000233"
$ws2.Cells.Item(33, 7).Value = 1

$ws2.Cells.Item(34, 1).Value = "Q6CDD3"
$ws2.Cells.Item(34, 2).Value = "YALI0C01411g"
$ws2.Cells.Item(34, 3).Value = "YALI0C01411g"
$ws2.Cells.Item(34, 4).Value = 51.419
# E34 carries no text, but keeps the quote-prefix style seen in the source file.
$ws2.Cells.Item(34, 5).Value = "'"
$ws2.Cells.Item(34, 5).Value = ""
$ws2.Cells.Item(34, 7).Value = 1

$ws2.Cells.Item(35, 1).Value = "Q6CAG1"

$ws2.Cells.Item(33, 6).Value = "Limits model when adding proteomics data. Calculatedby multiplying highest kcat/Km * Km (EC 1.14.19.41)"
$ws2.Cells.Item(34, 6).Value = "Limits model when adding proteomics data. Calculated from the specific activity of N. crassa (Brenda)"

$ws2.Cells.Item(35, 2).Value = "YALI0D03069g"
$ws2.Cells.Item(35, 3).Value = "YALI0D03069g"
$ws2.Cells.Item(35, 4).Value = 13.5
$ws2.Cells.Item(35, 7).Value = 1

$ws2.Cells.Item(36, 1).Value = "Q6CE88"
$ws2.Cells.Item(36, 2).Value = "YALI0B17644g"
$ws2.Cells.Item(36, 3).Value = "YALI0B17644g"
$ws2.Cells.Item(36, 4).Value = 14.3186
$ws2.Cells.Item(36, 5).Value = "y000236, y000237"
$ws2.Cells.Item(36, 7).Value = 1

$ws2.Cells.Item(35, 6).Value = "Limits model when adding proteomics data. Using kcat of E. coli (EC 2.1.2.2)"
$ws2.Cells.Item(36, 6).Value = "Limits model when adding proteomics data. Using standardKcat."

$ws2.Cells.Item(37, 1).Value = "Q6C5R8"
$ws2.Cells.Item(37, 2).Value = "YALI0E15730g"
$ws2.Cells.Item(37, 3).Value = "YALI0E15730g"
$ws2.Cells.Item(37, 4).Value = 1.43186
$ws2.Cells.Item(37, 7).Value = 1

$ws2.Cells.Item(37, 6).Value = "Limits model when adding proteomics data. Using standardKcat/10 instead."

# --- Widen column F (notes) on the "F 0.3" sheet for readability ---
$ws1.Columns.Item(6).ColumnWidth = 79.3

# --- Update the remembered cell selections on both sheets ---
$ws1.Activate()
[void]$ws1.Range("A46:E46").Select()
$ws2.Activate()
[void]$ws2.Range("F38").Select()
